# Renamed gim_damage_wood_step to gim_rolling_saw.
#
# The gimmick "gim_damage_wood_step" (Gimmick sheet, row 6 / id 2) is renamed
# to "gim_rolling_saw". Excel will drop the old shared string (it is no
# longer referenced anywhere) and append the new one to the shared string
# table, which also renumbers any shared strings that came after it.

$wb = $excel.ActiveWorkbook

$stage = $wb.Worksheets.Item("Stage")
$gimmick = $wb.Worksheets.Item("Gimmick")

# Update the gimmick's name.
$gimmick.Range("B6").Value = "gim_rolling_saw"

# Mirror the editor state (active sheet / selection) captured alongside the
# rename.
[void]$stage.Range("E5").Select()
[void]$gimmick.Activate()
[void]$gimmick.Range("E14").Select()
